# Update count_complexity. Update experiments according to the DISCLAIMER note
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2..5: rename test queries and update a couple of other values ---
$ws.Range("B2").Value = "manual_test_join_1"
$ws.Range("D2").Value = 170
$ws.Range("E2").Value = 1.5

$ws.Range("B3").Value = "high_level_test_join_1"
$ws.Range("E3").Value = 1

$ws.Range("B4").Value = "manual_test_join_2"
$ws.Range("D4").Value = 146
$ws.Range("E4").Value = 1.5

$ws.Range("B5").Value = "high_level_test_join_2"
$ws.Range("E5").Value = 1.5

# --- Recomputed "Level of Complexity (LoC)" values for the remaining rows ---
$loc = @{
    6  = 2
    7  = 1.5
    8  = 3
    9  = 2.5
    10 = 6.5
    11 = 1
    12 = 7
    13 = 1.5
    14 = 3.5
    16 = 4
    18 = 4
    19 = 1.5
    20 = 7.5
    21 = 1
    22 = 5
    23 = 2.5
    24 = 9
    25 = 2.5
    26 = 7.5
    27 = 1
    28 = 8
    29 = 1.5
    30 = 4
    31 = 1.5
    32 = 3.5
    33 = 2
}

foreach ($row in $loc.Keys) {
    $ws.Cells.Item($row, 5).Value = $loc[$row]
}
